# Roiger NEON Intern Methods Section - add fourth vignette about abundance
#
# This script reproduces, via Word COM-interop calls, the hand edit that:
#   1. Re-keys the existing inline comment from id 1 -> id 0 (a side effect
#      of Word's internal comment id bookkeeping - deleting the sole comment
#      and re-adding it at the same range makes it come back as id 0).
#   2. Replaces the placeholder "IX. Vignette X: TBD" heading with the real
#      title "IX. Vignette 4: Trends in abundance over time ", which in turn
#      moves Word's invisible "_GoBack" last-edit bookmark from the end of
#      the intro paragraph to the spot where the new text was typed.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Comment re-numbering: delete + recreate the single existing comment
#    on the same range, preserving its author/initials/text.
# ---------------------------------------------------------------------
$comment = $d.Comments.Item(1)
$commentScope = $comment.Scope
$commentText = $comment.Range.Text
$commentAuthor = $comment.Author
$commentInitials = $comment.Initial

$comment.Delete()
$newComment = $d.Comments.Add($commentScope, $commentText)
$newComment.Author = $commentAuthor
$newComment.Initial = $commentInitials

# ---------------------------------------------------------------------
# 2) Rewrite "Vignette X: TBD" -> "Vignette 4: Trends in abundance over
#    time ", splitting the run the way Word does when text is typed in
#    the middle of it, and dropping the _GoBack bookmark at the point
#    of the edit (which also removes it from its old location, since a
#    document can only have one _GoBack bookmark).
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute(" Vignette X: TBD", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$vignetteStart = $target.Start

# Split the run after " Vi" (offset 3) by dropping a throw-away bookmark
# there and removing it again - Word does not re-merge runs once split.
$splitPoint = $d.Range($vignetteStart + 3, $vignetteStart + 3)
$d.Bookmarks.Add("ZZZTempSplit", $splitPoint)
$d.Bookmarks.Item("ZZZTempSplit").Delete()

# Drop the _GoBack bookmark after "gnette " (offset 10) - this both
# splits the run there and relocates _GoBack from the intro paragraph.
$goBackPoint = $d.Range($vignetteStart + 10, $vignetteStart + 10)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# Replace "X: TBD" (offset 10-16) with the real vignette title.
$oldTitle = $d.Range($vignetteStart + 10, $vignetteStart + 16)
$oldTitle.Text = "4: Trends in abundance over time "
